$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row - copy the existing header style (from H1) onto the new header cells
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for columns I (I0) and J (IF)
$iValues = @(6, 8, 8, 8, 6, 9, 8, 9, 9, 9, 8, 6, 9, 2, 5)
$jValues = @(7, 8, 8, 8, 7, 9, 8, 9, 9, 9, 8, 6, 9, 2, 5)

for ($r = 0; $r -lt $iValues.Length; $r++) {
    $row = 2 + $r
    $ws.Cells.Item($row, 9).Value = $iValues[$r]
    $ws.Cells.Item($row, 10).Value = $jValues[$r]
}
